$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be written as literal text, independent of
# Excel's automatic number/date inference, without leaving any lasting
# NumberFormat/style change behind (matches source workbook where every
# D/E cell uses the default "General" style).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.911.27"
Set-TextValue $ws.Range("E2") "  -0.81%  "

Set-TextValue $ws.Range("D3") "2.213.88"
Set-TextValue $ws.Range("E3") "  -1.32%  "

Set-TextValue $ws.Range("E4") "  -0.05%  "

Set-TextValue $ws.Range("D5") "240.75"
Set-TextValue $ws.Range("E5") "  -2.51%  "

Set-TextValue $ws.Range("D6") "0.624"
Set-TextValue $ws.Range("E6") "  -0.80%  "

Set-TextValue $ws.Range("D7") "73.00"
Set-TextValue $ws.Range("E7") "  -2.14%  "

Set-TextValue $ws.Range("E8") "  +0.06%  "

Set-TextValue $ws.Range("D9") "0.603"
Set-TextValue $ws.Range("E9") "  -2.52%  "

Set-TextValue $ws.Range("D10") "42.68"
Set-TextValue $ws.Range("E10") "  +0.88%  "

Set-TextValue $ws.Range("D11") "0.0950"
Set-TextValue $ws.Range("E11") "  +0.69%  "

Set-TextValue $ws.Range("D12") "7.06"
Set-TextValue $ws.Range("E12") "  -1.57%  "

Set-TextValue $ws.Range("E13") "  -0.44%  "

Set-TextValue $ws.Range("D14") "2.547.08"
Set-TextValue $ws.Range("E14") "  -1.29%  "

Set-TextValue $ws.Range("D15") "14.19"
Set-TextValue $ws.Range("E15") "  -2.29%  "

Set-TextValue $ws.Range("D16") "0.833"
Set-TextValue $ws.Range("E16") "  -2.26%  "

Set-TextValue $ws.Range("D17") "2.223.27"
Set-TextValue $ws.Range("E17") "  -2.25%  "

Set-TextValue $ws.Range("D18") "41.775.98"
Set-TextValue $ws.Range("E18") "  -0.77%  "

Set-TextValue $ws.Range("E19") "  +9.82%  "

Set-TextValue $ws.Range("D20") "72.85"
Set-TextValue $ws.Range("E20") "  +0.89%  "

Set-TextValue $ws.Range("D21") "6.13"
Set-TextValue $ws.Range("E21") "  -0.03%  "

Set-TextValue $ws.Range("D22") "10.27"
Set-TextValue $ws.Range("E22") "  +15.38%  "

Set-TextValue $ws.Range("D23") "229.09"
Set-TextValue $ws.Range("E23") "  -1.13%  "

Set-TextValue $ws.Range("D24") "2.08"
Set-TextValue $ws.Range("E24") "  -6.02%  "

Set-TextValue $ws.Range("D25") "11.66"
Set-TextValue $ws.Range("E25") "  +2.06%  "

Set-TextValue $ws.Range("E26") "  -0.06%  "

Set-TextValue $ws.Range("E27") "  -0.85%  "

Set-TextValue $ws.Range("D28") "2.26"
Set-TextValue $ws.Range("E28") "  -2.13%  "

Set-TextValue $ws.Range("E29") "  +0.91%  "

Set-TextValue $ws.Range("D30") "166.81"
Set-TextValue $ws.Range("E30") "  -1.55%  "

Set-TextValue $ws.Range("D31") "20.54"
Set-TextValue $ws.Range("E31") "  -0.48%  "

Set-TextValue $ws.Range("D32") "5.62"
Set-TextValue $ws.Range("E32") "  +8.24%  "

Set-TextValue $ws.Range("D33") "0.0792"
Set-TextValue $ws.Range("E33") "  -3.63%  "

Set-TextValue $ws.Range("D34") "0.124"
Set-TextValue $ws.Range("E34") "  -0.65%  "

Set-TextValue $ws.Range("D35") "29.12"
Set-TextValue $ws.Range("E35") "  -7.22%  "

Set-TextValue $ws.Range("D36") "0.110"
Set-TextValue $ws.Range("E36") "  -8.52%  "

Set-TextValue $ws.Range("D37") "4.26"
Set-TextValue $ws.Range("E37") "  -5.02%  "

Set-TextValue $ws.Range("E38") "  -4.50%  "

Set-TextValue $ws.Range("D39") "13.50"
Set-TextValue $ws.Range("E39") "  -2.01%  "

Set-TextValue $ws.Range("D40") "65.51"
Set-TextValue $ws.Range("E40") "  +5.03%  "

Set-TextValue $ws.Range("E41") "  -3.33%  "

Set-TextValue $ws.Range("D42") "5.60"
Set-TextValue $ws.Range("E42") "  -2.70%  "

Set-TextValue $ws.Range("E43") "  -3.45%  "

Set-TextValue $ws.Range("D44") "8.64"
Set-TextValue $ws.Range("E44") "  -0.44%  "

Set-TextValue $ws.Range("D45") "103.59"
Set-TextValue $ws.Range("E45") "  -2.67%  "

Set-TextValue $ws.Range("E46") "  -2.21%  "

Set-TextValue $ws.Range("D47") "2.38"
Set-TextValue $ws.Range("E47") "  +3.75%  "

Set-TextValue $ws.Range("E48") "  -0.46%  "

Set-TextValue $ws.Range("E49") "  -0.83%  "

Set-TextValue $ws.Range("E50") "  -0.36%  "

Set-TextValue $ws.Range("D51") "2.419.37"
Set-TextValue $ws.Range("E51") "  -1.40%  "
